## Splits the "{first_name} {last_name}" placeholder run into separate
## "{firstName}", "{middleName}" and "{lastName}" runs (adding a middle
## initial placeholder), and relocates Word's "_GoBack" last-edit bookmark
## to sit between the middle- and last-name pieces, exactly as Word itself
## would do while a user is typing in that spot.

$d = $word.ActiveDocument

# Locate the placeholder text robustly (rather than assume a paragraph index).
$oldText = "{first_name} {last_name}"
$newText = "{firstName} {middleName} {lastName}"

$content = $d.Content
$allText = $content.Text
$idx = $allText.IndexOf($oldText)

$nameStart = $content.Start + $idx
$nameEnd   = $nameStart + $oldText.Length

# --- Step 1: detach the run's leading run (e.g. the preceding <w:tab/>)
# from the text we are about to rewrite, so it keeps its own run instead of
# being absorbed into the text run. Dropping a bookmark at the boundary and
# immediately deleting it splits the underlying run without altering it.
$splitPoint = $d.Range($nameStart, $nameStart)
$d.Bookmarks.Add("_zz_tmp_split_0", $splitPoint)
$d.Bookmarks("_zz_tmp_split_0").Delete()

# --- Step 2: rewrite the placeholder text in place.
$nameRange = $d.Range($nameStart, $nameEnd)
$nameRange.Text = $newText

# --- Step 3: split the new text into the individual placeholder runs.
# Offsets within $newText ("{firstName} {middleName} {lastName}"):
#   {firstName}        -> [0,11)
#   " {middleName}"     -> [11,24)
#   " {lastN"            -> [24,31)
#   "ame}"                -> [31,35)
$split3 = $nameStart + 31   # between " {lastN" and "ame}"
$split1 = $nameStart + 11   # between "{firstName}" and " {middleName}"
$split2 = $nameStart + 24   # between " {middleName}" and " {lastN" (bookmark goes here)

# Split off the tail piece first (while the run has no leading space yet),
# so the resulting fragments don't inherit an unnecessary xml:space="preserve".
$p3 = $d.Range($split3, $split3)
$d.Bookmarks.Add("_zz_tmp_split_3", $p3)
$d.Bookmarks("_zz_tmp_split_3").Delete()

$p1 = $d.Range($split1, $split1)
$d.Bookmarks.Add("_zz_tmp_split_1", $p1)
$d.Bookmarks("_zz_tmp_split_1").Delete()

# --- Step 4: drop Word's "_GoBack" bookmark at the new edit location
# (between the middle-name and last-name pieces). Word keeps only one
# "_GoBack" bookmark in the document, so adding it here automatically
# removes it from its previous location (at the end of the
# researchInterests paragraph).
$p2 = $d.Range($split2, $split2)
$d.Bookmarks.Add("_GoBack", $p2)
